$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new "Position" table
$ws.Range("A16").Value = "Position #"
$ws.Range("B16").Value = "Position xval (cm)"

# Data rows 17-23
$data = @(
    @(10, 17.8, 0.1804, 0.0985),
    @(8,  24.4, 0.2277, 0.1501),
    @(6,  32.2, 0.2854, 0.1973),
    @(5,  36.6, 0.3195, 0.2321),
    @(4,  40.8, 0.352, 0.2617),
    @(2,  50, 0.4264, 0.3273),
    @(1,  55, 0.4684, 0.3536)
)

$row = 17
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Column E formulas: row 17 is a standalone formula, rows 18-23 form a
# shared-formula group (assigning to the whole range at once is what makes
# the engine emit a shared formula group rather than per-cell formulas).
$ws.Range("E17").Formula = "=C17-D17"
$ws.Range("E18:E23").Formula = "=C18-D18"

$ws.Range("F22").Select()
